# Updates the "cryptos" price/volume table (columns D and E, rows 2-51)
# with refreshed values from the latest GitHub Actions scrape run.
#
# Column D ("Price") holds plain text in this sheet (prices use '.' as a
# thousands separator, e.g. "62.812.76"), so most new values already fail
# Excel's automatic number parsing and stay text on their own. A handful of
# new prices (e.g. "1.00", "577.75") *do* look like valid numbers to Excel,
# so for those we briefly force the cell to Text format before assigning
# the value, then restore the cell's original ("Normal") style so no
# stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.812.76'
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").Value = '3.461.07'
$ws.Range("E3").Value = '  +2.04%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '577.75'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '147.60'
$c.Style = "Normal"
$ws.Range("D7").Value = '3.460.98'
$ws.Range("E7").Value = '  +2.03%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +1.73%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.69'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("E12").Value = '  +4.14%  '
$ws.Range("D13").Value = '4.053.35'
$ws.Range("E13").Value = '  +2.02%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '29.74'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +6.43%  '
$ws.Range("E15").Value = '  +2.85%  '
$ws.Range("D16").Value = '3.464.59'
$ws.Range("E16").Value = '  +1.72%  '
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '62.814.96'
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("E19").Value = '  +3.55%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.28'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +5.18%  '
$ws.Range("E21").Value = '  +1.35%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '388.56'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.31%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.555'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '3.603.18'
$ws.Range("E26").Value = '  +2.02%  '
$ws.Range("E27").Value = '  +1.68%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.179'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -6.54%  '
$ws.Range("E29").Value = '  +1.78%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.13%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.12'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.03%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.12'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  -2.25%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '23.62'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("E36").Value = '  +2.04%  '
$ws.Range("E37").Value = '  +3.74%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '31.73'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +19.76%  '
$ws.Range("E39").Value = '  +6.97%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '169.68'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.78%  '
$ws.Range("D41").Value = '3.499.87'
$ws.Range("E41").Value = '  +2.10%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0752'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.27%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.797'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.16%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '42.35'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.25%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '4.45'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("E46").Value = '  +3.28%  '
$ws.Range("E47").Value = '  +3.88%  '
$ws.Range("D48").Value = '2.598.17'
$ws.Range("E48").Value = '  +5.70%  '
$ws.Range("E49").Value = '  +11.78%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '22.90'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.38%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.71'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.87%  '
